$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column B (was "ANALOGOUS DRUG COMPOUNDS" -> "ANALOGOUS DRUG COMPOUND")
$ws.Range("B1").Value = "ANALOGOUS DRUG COMPOUND"

# Data rows: precursor drug compound (col A) and analogous drug compound (col B)
$data = @(
    @("DEGARELIX", "ACYLINE"),
    @("DEGARELIX", "SATEREOTIDE TETRAXETAN"),
    @("DEGARELIX", "CETRORELIX"),
    @("DEGARELIX", "GANIRELIX"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "DESLORELIN"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "NAFARELIN"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "GONADORELIN"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "AFAMELANOTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "ALSACTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "CZEN 202"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "MODIMELANOTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "TASPOGLUTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "TETRACOSACTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "BREMELANOTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "COTATDUTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "TRIDECTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "AMY-101"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "TRY-120027"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "ANGIOTENSINAMIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "SEMAGLUTIDE "),
    @("ELIGARD (LEUPROLIDE ACETATE)", "ALBUVIRTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "ANGIOTENSIN 1-7"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "LIXISENATIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "SETMELANOTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "SARALASIN"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "PRAMLINTIDE"),
    @("ELIGARD (LEUPROLIDE ACETATE)", "MUREPAVADIN"),
    @("FLUTAMIDE", "(R)-3-BROMO-2-HYDROXY-2-METHYL-N-[4-NITRO-3-(TRIFLUOROMETHYL)PHENYL]PROPANAMIDE AND CANCER"),
    @("GOSERELIN ACETATE", "ANGIOTENSIN II"),
    @("GOSERELIN ACETATE", "ACLERASTIDE")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Adjust column widths to match target (values chosen so the stored OOXML
# width quantizes to 13 and 20.5 respectively under this engine's rounding)
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth = 19.666666666666668

# Selection as recorded in diff
$ws.Range("B21").Select()
